$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Strategy:
#  Word's InsertParagraphBefore()/After() copies the pPr of the
#  paragraph it is attached to, which for our numbered "question"
#  paragraphs means the new paragraph would wrongly inherit the
#  Prrafodelista/numPr list formatting. However the LAST paragraph in
#  the original document (the "Explicacion..." one) already has the
#  plain formatting we want for all the new answer paragraphs
#  (just <w:jc w:val="both"/>, no pStyle/numPr). Copy/Paste in this
#  engine preserves the copied paragraph's own formatting regardless
#  of where it is pasted, so we build one clean "stamp" paragraph and
#  Copy+Paste it everywhere we need a new plain paragraph, then set
#  its text.
# ------------------------------------------------------------------

$lastIdx = $d.Paragraphs.Count
$pLast = $d.Paragraphs($lastIdx)
$pLast.Range.InsertParagraphAfter()
$stamp = $d.Paragraphs($lastIdx + 1)
$stamp.Range.Text = "STAMP"
$stamp.Range.Copy()

function New-CleanParagraphBefore($paraIndex, [string]$text) {
    $target = $d.Paragraphs($paraIndex)
    $rng = $d.Range($target.Range.Start, $target.Range.Start)
    $rng.Paste()
    $newPara = $d.Paragraphs($paraIndex)
    $newPara.Range.Text = $text
}

# Paragraphs 3..8 are, at this point, the six interview questions
# (Q1..Q6); paragraph 9 is still the original "Explicacion..." one.
# Insert the plain answer paragraphs after Q1..Q5 (i.e. before the
# following question), working from the bottom up so indices stay
# valid.

New-CleanParagraphBefore 8 "No se brinda capacitación dentro de Grupo Promesa, solamente a operarios como requisito de Dirección Nacional de Medicamentos, para el manejo de farmacia y su renovación. "
New-CleanParagraphBefore 7 "No se cuentan con manuales administrativos."
New-CleanParagraphBefore 6 "No existe una conexión entre la ubicación del hospital y la ubicación de la clínica y la farmacia, pero si existe entre las computadoras de cada infraestructura."
New-CleanParagraphBefore 5 "Todas las áreas cuentan con computadora a excepción del área de enfermería."
New-CleanParagraphBefore 4 "No existe una frecuencia fija para la actualización del equipo, cada vez que algún equipo ya no funcione es cambiado y según surjan las necesidades se va adquiriendo nuevo equipo."

# Now paragraph order is:
#  1 Title
#  2 Subtitle
#  3 Q1
#  4 answer1
#  5 Q2
#  6 answer2
#  7 Q3
#  8 answer3
#  9 Q4
# 10 answer4
# 11 Q5
# 12 answer5
# 13 Q6
# 14 "Explicacion..." (original paragraph, still with mid-text bookmark)
# 15 STAMP

# Insert the two new plain paragraphs before the original
# "Explicacion..." paragraph (now index 14).
New-CleanParagraphBefore 14 "No se cuenta con un estándar los códigos son establecidos al azar."
New-CleanParagraphBefore 15 "Explicación sobre la dificultad para migrar datos de una aplicación a otra debido al mal uso de estándares y formatos de datos para su registro."

# Paragraph order now:
#  ...
# 13 Q6
# 14 "No se cuenta con un estándar..."
# 15 "Explicación sobre la dificultad..."
# 16 original "Explicacion..." paragraph (2 runs + bookmark in the middle)
# 17 STAMP

$origIdx = 16
$origPara = $d.Paragraphs($origIdx)
$origStart = $origPara.Range.Start

$run1Old = "Explicación sobre la dificultad para migrar datos de una aplic"
$len1 = $run1Old.Length

$newRun1Text = "Debido a que existe una incorrecta manipulación de los registros y la falta de "
$newRun2Text = "normalización de las tablas no es posible la migración y se opta por un ingreso manual de los datos."

# Edit the first run's text in place (this does not disturb the
# bookmark that immediately follows it).
$r1 = $d.Range($origStart, $origStart + $len1)
$r1.Text = $newRun1Text

# Recompute the boundary right after run1 now that its text changed -
# the bookmark sits exactly there, at the run1/run2 split point.
$newRun1End = $origStart + $newRun1Text.Length
$currentParaEnd = $d.Paragraphs($origIdx).Range.End

# Edit the second run's text in place too (again without touching the
# bookmark, which sits right before this range).
$r2 = $d.Range($newRun1End, $currentParaEnd - 1)
$r2.Text = $newRun2Text

# ------------------------------------------------------------------
# Move the _GoBack bookmark (currently sitting between run1 and run2)
# to the very end of the paragraph, after run2, right before the
# paragraph mark. Bookmarks.Add() with a truly empty/collapsed range
# at certain end-of-document-ish offsets snaps to the start of the
# document in this engine, so instead we temporarily insert a marker
# character at the destination, bookmark the (non-empty) range that
# wraps it, then delete just that marker character - the bookmark
# collapses naturally to the correct spot and survives the in-place
# text edit (as already demonstrated for run edits above).
# ------------------------------------------------------------------

$finalPara = $d.Paragraphs($origIdx)
$finalEnd = $finalPara.Range.End
$markerPos = $d.Range($finalEnd - 1, $finalEnd - 1)
$markerPos.InsertAfter("X")

$finalEnd2 = $d.Paragraphs($origIdx).Range.End
$markerRange = $d.Range($finalEnd2 - 2, $finalEnd2 - 1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$finalEnd3 = $d.Paragraphs($origIdx).Range.End
$markerRange2 = $d.Range($finalEnd3 - 2, $finalEnd3 - 1)
$markerRange2.Text = ""

# ------------------------------------------------------------------
# Clean up: remove the STAMP paragraph we used as a copy/paste source.
# ------------------------------------------------------------------

$stampIdx = $d.Paragraphs.Count
$stampPara = $d.Paragraphs($stampIdx)
$stampPara.Range.Delete()

Write-Host "Final paragraph count:" $d.Paragraphs.Count
